$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet needs to hold the *text* "1" (it keeps its
# existing shared-string type/style - only the displayed text changes).
#
# Assigning a plain numeric-looking string straight to .Value (e.g. "1")
# makes Excel store it as a real number, which would change the cell's
# type away from a shared string. Using a leading apostrophe forces text,
# but doing that directly on B11 also stamps a brand-new "quote prefixed"
# style onto the cell, which would needlessly change its style id.
#
# So: build the text value "1" in a scratch cell far outside the used
# range, then copy/paste just the *value* into B11. Paste-values-only
# carries the text over without touching B11's existing formatting/style.
$helper = $ws.Range("Z1")
$helper.Value = "'1"

$helper.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

# Tidy up: clear the scratch cell and the marching-ants clipboard state.
$helper.Clear()
$excel.CutCopyMode = $false
